$d = $word.ActiveDocument

$pairs = @(
    @("91×76=6916", "67×78=5226"),
    @("31×61=1891", "98×87=8526"),
    @("56×46=2576", "58×73=4234"),
    @("62×83=5146", "44×16=704"),
    @("93×12=1116", "84×74=6216"),
    @("46×39=1794", "23×99=2277"),
    @("54×26=1404", "51×71=3621"),
    @("65×77=5005", "77×87=6699"),
    @("80×21=1680", "62×12=744"),
    @("93×15=1395", "65×39=2535"),
    @("36×76=2736", "81×25=2025"),
    @("65×76=4940", "85×70=5950"),
    @("13×91=1183", "31×13=403"),
    @("76×62=4712", "99×62=6138"),
    @("30×99=2970", "31×68=2108"),
    @("55×97=5335", "88×45=3960"),
    @("47×21=987", "88×73=6424"),
    @("12×99=1188", "16×57=912"),
    @("36×22=792", "12×21=252"),
    @("65×32=2080", "23×92=2116"),
    @("91×62=5642", "75×19=1425"),
    @("53×67=3551", "67×24=1608"),
    @("83×55=4565", "95×84=7980"),
    @("74×90=6660", "33×51=1683"),
    @("55×56=3080", "17×87=1479")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
